# Adds the "resourceDeleteConfirmation" notification setting + its two
# NOTIFY_TEMPLATE rows (email + pm), per commit "[ADD] add resource delete template".

$wb = $excel.ActiveWorkbook

# --- NOTIFY_SEND_SETTING: new row 21 ---------------------------------------
$wsSetting = $wb.Worksheets.Item("NOTIFY_SEND_SETTING")

$wsSetting.Range("E21").Value = 'resourceDeleteConfirmation'
$wsSetting.Range("F21").Value = 'resourceDeleteConfirmation'
$wsSetting.Range("G21").Value = '资源删除确认通知'
$wsSetting.Range("H21").Value = '删除Devops资源时，给相关用户发送验证码'
$wsSetting.Range("I21").Value = "site"
$wsSetting.Range("J21").Value = 0
$wsSetting.Range("K21").Value = 1
$wsSetting.Range("L21").Value = 1
$wsSetting.Range("M21").Value = 1

# --- NOTIFY_TEMPLATE: new rows 21 (email) and 22 (pm) -----------------------
$wsTemplate = $wb.Worksheets.Item("NOTIFY_TEMPLATE")

$wsTemplate.Range("E21").Value = 'resourceDeleteConfirmation-preset'
$wsTemplate.Range("F21").Value = 'resourceDeleteConfirmation-email'
$wsTemplate.Range("G21").Value = '资源删除确认通知'
$wsTemplate.Range("H21").Value = "email"
$wsTemplate.Range("I21").Value = 1
$wsTemplate.Range("J21").Value = 'resourceDeleteConfirmation'
$wsTemplate.Range("K21").Value = 'Choerodon验证邮件'
$wsTemplate.Range("N21").Value = '您好，${user}正在${env}环境下执行删除${object}"${objectName}"的操作，验证码为：${verificationCode}；确认后，需将此验证码提供给操作者${user}完成删除操作。验证码10分钟内有效。'

$wsTemplate.Range("E22").Value = 'resourceDeleteConfirmation-preset'
$wsTemplate.Range("F22").Value = 'resourceDeleteConfirmation-pm'
$wsTemplate.Range("G22").Value = '资源删除确认通知'
$wsTemplate.Range("H22").Value = "pm"
$wsTemplate.Range("I22").Value = 1
$wsTemplate.Range("J22").Value = 'resourceDeleteConfirmation'
$wsTemplate.Range("L22").Value = '删除操作验证码'
$wsTemplate.Range("M22").Value = '${user}正在${env}环境下执行删除${object}"${objectName}"的操作，验证码为：${verificationCode}；确认后，需将此验证码提供给操作者${user}完成删除操作。验证码10分钟内有效。'

# --- view state: mirror the saved workbook (active tab + selections) -------
$wsTemplate.Cells.Select()

$wsSetting.Activate()
$wsSetting.Range("D7:M21").Select()
